$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 452, shifting rows 452:549 down to 453:550
$ws.Rows.Item(452).Insert()

# Populate the newly inserted row 452 with the new record's data
$ws.Cells.Item(452, 1).Value = 3
$ws.Cells.Item(452, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(452, 3).Value = "Coquimbo"
$ws.Cells.Item(452, 4).Value = 45244
$ws.Cells.Item(452, 5).Value = 5
$ws.Cells.Item(452, 6).Value = 100112001
$ws.Cells.Item(452, 7).Value = "Berenjena"
$ws.Cells.Item(452, 8).Value = "Sin especificar"
$ws.Cells.Item(452, 9).Value = "Primera"
$ws.Cells.Item(452, 10).Value = 50
$ws.Cells.Item(452, 11).Value = 9000
$ws.Cells.Item(452, 12).Value = 9000
$ws.Cells.Item(452, 13).Value = 9000
$ws.Cells.Item(452, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(452, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(452, 16).Value = 150
$ws.Cells.Item(452, 17).Value = 60
$ws.Cells.Item(452, 18).Value = "Hortaliza"
